# Auto-generated edit script: update stats values in rows 4-13 per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 0.289
$ws.Cells.Item(4, 5).Value = 0.168
$ws.Cells.Item(4, 6).Value = 0.028
$ws.Cells.Item(4, 7).Value = 0.168
$ws.Cells.Item(4, 8).Value = 0.2
$ws.Cells.Item(4, 9).Value = 0.026
$ws.Cells.Item(4, 10).Value = 0.161
$ws.Cells.Item(4, 11).Value = 0.307
$ws.Cells.Item(4, 12).Value = 0.099
$ws.Cells.Item(4, 13).Value = 0.314
$ws.Cells.Item(4, 14).Value = 0.254
$ws.Cells.Item(4, 16).Value = 0.149
$ws.Cells.Item(4, 17).Value = 0.478
$ws.Cells.Item(4, 18).Value = 0.222
$ws.Cells.Item(4, 20).Value = 0.244
$ws.Cells.Item(4, 22).Value = 0.289
$ws.Cells.Item(4, 23).Value = 0.242
$ws.Cells.Item(4, 24).Value = 0.043
$ws.Cells.Item(4, 26).Value = 0.431
$ws.Cells.Item(4, 27).Value = 0.13
$ws.Cells.Item(4, 28).Value = 0.361
$ws.Cells.Item(4, 29).Value = 0.117
$ws.Cells.Item(4, 31).Value = 0.082
$ws.Cells.Item(4, 32).Value = 0.713
$ws.Cells.Item(4, 34).Value = 0.325
$ws.Cells.Item(4, 35).Value = 0.656
$ws.Cells.Item(4, 36).Value = 0.165
$ws.Cells.Item(4, 37).Value = 0.406
$ws.Cells.Item(4, 38).Value = 0.671
$ws.Cells.Item(4, 41).Value = 0.68
$ws.Cells.Item(5, 2).Value = 0.822
$ws.Cells.Item(5, 3).Value = 0.146
$ws.Cells.Item(5, 4).Value = 0.382
$ws.Cells.Item(5, 5).Value = 0.667
$ws.Cells.Item(5, 6).Value = 0.222
$ws.Cells.Item(5, 7).Value = 0.471
$ws.Cells.Item(5, 8).Value = 0.8
$ws.Cells.Item(5, 9).Value = 0.16
$ws.Cells.Item(5, 10).Value = 0.4
$ws.Cells.Item(5, 11).Value = 0.6
$ws.Cells.Item(5, 12).Value = 0.24
$ws.Cells.Item(5, 13).Value = 0.49
$ws.Cells.Item(5, 14).Value = 0.8
$ws.Cells.Item(5, 15).Value = 0.16
$ws.Cells.Item(5, 16).Value = 0.4
$ws.Cells.Item(5, 17).Value = 0.533
$ws.Cells.Item(5, 18).Value = 0.249
$ws.Cells.Item(5, 19).Value = 0.499
$ws.Cells.Item(5, 20).Value = 0.511
$ws.Cells.Item(5, 21).Value = 0.25
$ws.Cells.Item(5, 22).Value = 0.5
$ws.Cells.Item(5, 23).Value = 0.733
$ws.Cells.Item(5, 24).Value = 0.196
$ws.Cells.Item(5, 25).Value = 0.442
$ws.Cells.Item(5, 26).Value = 0.8
$ws.Cells.Item(5, 27).Value = 0.16
$ws.Cells.Item(5, 28).Value = 0.4
$ws.Cells.Item(5, 29).Value = 0.711
$ws.Cells.Item(5, 30).Value = 0.205
$ws.Cells.Item(5, 31).Value = 0.453
$ws.Cells.Item(5, 32).Value = 0.956
$ws.Cells.Item(5, 33).Value = 0.042
$ws.Cells.Item(5, 34).Value = 0.206
$ws.Cells.Item(5, 35).Value = 0.778
$ws.Cells.Item(5, 36).Value = 0.173
$ws.Cells.Item(5, 37).Value = 0.416
$ws.Cells.Item(5, 38).Value = 0.911
$ws.Cells.Item(5, 39).Value = 0.081
$ws.Cells.Item(5, 40).Value = 0.285
$ws.Cells.Item(5, 41).Value = 0.882
$ws.Cells.Item(6, 2).Value = 0.428
$ws.Cells.Item(6, 5).Value = 0.268
$ws.Cells.Item(6, 8).Value = 0.32
$ws.Cells.Item(6, 11).Value = 0.406
$ws.Cells.Item(6, 14).Value = 0.386
$ws.Cells.Item(6, 17).Value = 0.504
$ws.Cells.Item(6, 20).Value = 0.33
$ws.Cells.Item(6, 23).Value = 0.364
$ws.Cells.Item(6, 26).Value = 0.5600000000000001
$ws.Cells.Item(6, 29).Value = 0.201
$ws.Cells.Item(6, 32).Value = 0.8169999999999999
$ws.Cells.Item(6, 35).Value = 0.712
$ws.Cells.Item(6, 38).Value = 0.773
$ws.Cells.Item(6, 41).Value = 0.767
$ws.Cells.Item(7, 2).Value = 0.601
$ws.Cells.Item(7, 5).Value = 0.418
$ws.Cells.Item(7, 8).Value = 0.5
$ws.Cells.Item(7, 11).Value = 0.504
$ws.Cells.Item(7, 14).Value = 0.5590000000000001
$ws.Cells.Item(7, 17).Value = 0.521
$ws.Cells.Item(7, 20).Value = 0.419
$ws.Cells.Item(7, 23).Value = 0.521
$ws.Cells.Item(7, 26).Value = 0.6830000000000001
$ws.Cells.Item(7, 29).Value = 0.353
$ws.Cells.Item(7, 32).Value = 0.895
$ws.Cells.Item(7, 35).Value = 0.75
$ws.Cells.Item(7, 38).Value = 0.85
$ws.Cells.Item(7, 41).Value = 0.832
$ws.Cells.Item(8, 2).Value = 0.742
$ws.Cells.Item(8, 3).Value = 0.149
$ws.Cells.Item(8, 4).Value = 0.386
$ws.Cells.Item(8, 5).Value = 0.5629999999999999
$ws.Cells.Item(8, 8).Value = 0.697
$ws.Cells.Item(8, 9).Value = 0.158
$ws.Cells.Item(8, 10).Value = 0.398
$ws.Cells.Item(8, 11).Value = 0.531
$ws.Cells.Item(8, 13).Value = 0.46
$ws.Cells.Item(8, 14).Value = 0.713
$ws.Cells.Item(8, 15).Value = 0.157
$ws.Cells.Item(8, 16).Value = 0.396
$ws.Cells.Item(8, 17).Value = 0.509
$ws.Cells.Item(8, 19).Value = 0.484
$ws.Cells.Item(8, 20).Value = 0.445
$ws.Cells.Item(8, 23).Value = 0.662
$ws.Cells.Item(8, 24).Value = 0.182
$ws.Cells.Item(8, 25).Value = 0.426
$ws.Cells.Item(8, 26).Value = 0.737
$ws.Cells.Item(8, 27).Value = 0.157
$ws.Cells.Item(8, 28).Value = 0.396
$ws.Cells.Item(8, 29).Value = 0.596
$ws.Cells.Item(8, 30).Value = 0.189
$ws.Cells.Item(8, 31).Value = 0.435
$ws.Cells.Item(8, 32).Value = 0.879
$ws.Cells.Item(8, 33).Value = 0.06
$ws.Cells.Item(8, 34).Value = 0.244
$ws.Cells.Item(8, 35).Value = 0.77
$ws.Cells.Item(8, 36).Value = 0.172
$ws.Cells.Item(8, 37).Value = 0.415
$ws.Cells.Item(8, 38).Value = 0.878
$ws.Cells.Item(8, 39).Value = 0.08599999999999999
$ws.Cells.Item(8, 40).Value = 0.294
$ws.Cells.Item(8, 41).Value = 0.842
$ws.Cells.Item(9, 2).Value = 0.644
$ws.Cells.Item(9, 3).Value = 0.229
$ws.Cells.Item(9, 4).Value = 0.479
$ws.Cells.Item(9, 5).Value = 0.444
$ws.Cells.Item(9, 6).Value = 0.247
$ws.Cells.Item(9, 7).Value = 0.497
$ws.Cells.Item(9, 8).Value = 0.578
$ws.Cells.Item(9, 9).Value = 0.244
$ws.Cells.Item(9, 10).Value = 0.494
$ws.Cells.Item(9, 11).Value = 0.444
$ws.Cells.Item(9, 12).Value = 0.247
$ws.Cells.Item(9, 13).Value = 0.497
$ws.Cells.Item(9, 14).Value = 0.6
$ws.Cells.Item(9, 15).Value = 0.24
$ws.Cells.Item(9, 16).Value = 0.49
$ws.Cells.Item(9, 17).Value = 0.467
$ws.Cells.Item(9, 20).Value = 0.356
$ws.Cells.Item(9, 21).Value = 0.229
$ws.Cells.Item(9, 22).Value = 0.479
$ws.Cells.Item(9, 23).Value = 0.556
$ws.Cells.Item(9, 24).Value = 0.247
$ws.Cells.Item(9, 25).Value = 0.497
$ws.Cells.Item(9, 26).Value = 0.644
$ws.Cells.Item(9, 27).Value = 0.229
$ws.Cells.Item(9, 28).Value = 0.479
$ws.Cells.Item(9, 29).Value = 0.489
$ws.Cells.Item(9, 30).Value = 0.25
$ws.Cells.Item(9, 31).Value = 0.5
$ws.Cells.Item(9, 32).Value = 0.756
$ws.Cells.Item(9, 33).Value = 0.185
$ws.Cells.Item(9, 34).Value = 0.43
$ws.Cells.Item(9, 35).Value = 0.756
$ws.Cells.Item(9, 36).Value = 0.185
$ws.Cells.Item(9, 37).Value = 0.43
$ws.Cells.Item(9, 38).Value = 0.822
$ws.Cells.Item(9, 39).Value = 0.146
$ws.Cells.Item(9, 40).Value = 0.382
$ws.Cells.Item(9, 41).Value = 0.778
$ws.Cells.Item(10, 2).Value = 0.778
$ws.Cells.Item(10, 3).Value = 0.173
$ws.Cells.Item(10, 4).Value = 0.416
$ws.Cells.Item(10, 5).Value = 0.6
$ws.Cells.Item(10, 6).Value = 0.24
$ws.Cells.Item(10, 7).Value = 0.49
$ws.Cells.Item(10, 8).Value = 0.733
$ws.Cells.Item(10, 9).Value = 0.196
$ws.Cells.Item(10, 10).Value = 0.442
$ws.Cells.Item(10, 11).Value = 0.6
$ws.Cells.Item(10, 12).Value = 0.24
$ws.Cells.Item(10, 13).Value = 0.49
$ws.Cells.Item(10, 14).Value = 0.778
$ws.Cells.Item(10, 15).Value = 0.173
$ws.Cells.Item(10, 16).Value = 0.416
$ws.Cells.Item(10, 17).Value = 0.533
$ws.Cells.Item(10, 18).Value = 0.249
$ws.Cells.Item(10, 19).Value = 0.499
$ws.Cells.Item(10, 20).Value = 0.511
$ws.Cells.Item(10, 21).Value = 0.25
$ws.Cells.Item(10, 22).Value = 0.5
$ws.Cells.Item(10, 23).Value = 0.733
$ws.Cells.Item(10, 24).Value = 0.196
$ws.Cells.Item(10, 25).Value = 0.442
$ws.Cells.Item(10, 26).Value = 0.8
$ws.Cells.Item(10, 27).Value = 0.16
$ws.Cells.Item(10, 28).Value = 0.4
$ws.Cells.Item(10, 29).Value = 0.6
$ws.Cells.Item(10, 30).Value = 0.24
$ws.Cells.Item(10, 31).Value = 0.49
$ws.Cells.Item(10, 32).Value = 0.956
$ws.Cells.Item(10, 33).Value = 0.042
$ws.Cells.Item(10, 34).Value = 0.206
$ws.Cells.Item(10, 35).Value = 0.778
$ws.Cells.Item(10, 36).Value = 0.173
$ws.Cells.Item(10, 37).Value = 0.416
$ws.Cells.Item(10, 38).Value = 0.911
$ws.Cells.Item(10, 39).Value = 0.081
$ws.Cells.Item(10, 40).Value = 0.285
$ws.Cells.Item(10, 41).Value = 0.882
$ws.Cells.Item(11, 2).Value = 0.822
$ws.Cells.Item(11, 3).Value = 0.146
$ws.Cells.Item(11, 4).Value = 0.382
$ws.Cells.Item(11, 5).Value = 0.667
$ws.Cells.Item(11, 6).Value = 0.222
$ws.Cells.Item(11, 7).Value = 0.471
$ws.Cells.Item(11, 8).Value = 0.8
$ws.Cells.Item(11, 9).Value = 0.16
$ws.Cells.Item(11, 10).Value = 0.4
$ws.Cells.Item(11, 11).Value = 0.6
$ws.Cells.Item(11, 12).Value = 0.24
$ws.Cells.Item(11, 13).Value = 0.49
$ws.Cells.Item(11, 14).Value = 0.8
$ws.Cells.Item(11, 15).Value = 0.16
$ws.Cells.Item(11, 16).Value = 0.4
$ws.Cells.Item(11, 17).Value = 0.533
$ws.Cells.Item(11, 18).Value = 0.249
$ws.Cells.Item(11, 19).Value = 0.499
$ws.Cells.Item(11, 20).Value = 0.511
$ws.Cells.Item(11, 21).Value = 0.25
$ws.Cells.Item(11, 22).Value = 0.5
$ws.Cells.Item(11, 23).Value = 0.733
$ws.Cells.Item(11, 24).Value = 0.196
$ws.Cells.Item(11, 25).Value = 0.442
$ws.Cells.Item(11, 26).Value = 0.8
$ws.Cells.Item(11, 27).Value = 0.16
$ws.Cells.Item(11, 28).Value = 0.4
$ws.Cells.Item(11, 29).Value = 0.644
$ws.Cells.Item(11, 30).Value = 0.229
$ws.Cells.Item(11, 31).Value = 0.479
$ws.Cells.Item(11, 32).Value = 0.956
$ws.Cells.Item(11, 33).Value = 0.042
$ws.Cells.Item(11, 34).Value = 0.206
$ws.Cells.Item(11, 35).Value = 0.778
$ws.Cells.Item(11, 36).Value = 0.173
$ws.Cells.Item(11, 37).Value = 0.416
$ws.Cells.Item(11, 38).Value = 0.911
$ws.Cells.Item(11, 39).Value = 0.081
$ws.Cells.Item(11, 40).Value = 0.285
$ws.Cells.Item(11, 41).Value = 0.882
$ws.Cells.Item(12, 2).Value = 1.378
$ws.Cells.Item(12, 3).Value = 0.668
$ws.Cells.Item(12, 4).Value = 0.8169999999999999
$ws.Cells.Item(12, 5).Value = 1.633
$ws.Cells.Item(12, 6).Value = 1.032
$ws.Cells.Item(12, 7).Value = 1.016
$ws.Cells.Item(12, 8).Value = 1.556
$ws.Cells.Item(12, 9).Value = 1.191
$ws.Cells.Item(12, 10).Value = 1.091
$ws.Cells.Item(12, 11).Value = 1.407
$ws.Cells.Item(12, 12).Value = 0.538
$ws.Cells.Item(12, 13).Value = 0.733
$ws.Cells.Item(12, 14).Value = 1.389
$ws.Cells.Item(12, 15).Value = 0.571
$ws.Cells.Item(12, 16).Value = 0.756
$ws.Cells.Item(12, 26).Value = 1.25
$ws.Cells.Item(12, 27).Value = 0.299
$ws.Cells.Item(12, 28).Value = 0.546
$ws.Cells.Item(12, 29).Value = 2
$ws.Cells.Item(12, 30).Value = 3.812
$ws.Cells.Item(12, 31).Value = 1.953
$ws.Cells.Item(12, 32).Value = 1.233
$ws.Cells.Item(12, 33).Value = 0.225
$ws.Cells.Item(12, 34).Value = 0.474
$ws.Cells.Item(12, 36).Value = 0.028
$ws.Cells.Item(12, 37).Value = 0.167
$ws.Cells.Item(12, 38).Value = 1.098
$ws.Cells.Item(12, 39).Value = 0.08799999999999999
$ws.Cells.Item(12, 40).Value = 0.297
$ws.Cells.Item(12, 41).Value = 1.12
$ws.Cells.Item(13, 2).Value = 3.533
$ws.Cells.Item(13, 3).Value = 1.404
$ws.Cells.Item(13, 4).Value = 1.185
$ws.Cells.Item(13, 5).Value = 4.564
$ws.Cells.Item(13, 6).Value = 0.707
$ws.Cells.Item(13, 7).Value = 0.841
$ws.Cells.Item(13, 8).Value = 4.524
$ws.Cells.Item(13, 9).Value = 0.916
$ws.Cells.Item(13, 10).Value = 0.957
$ws.Cells.Item(13, 11).Value = 2.3
$ws.Cells.Item(13, 12).Value = 0.61
$ws.Cells.Item(13, 13).Value = 0.781
$ws.Cells.Item(13, 14).Value = 3.333
$ws.Cells.Item(13, 15).Value = 0.756
$ws.Cells.Item(13, 16).Value = 0.869
$ws.Cells.Item(13, 26).Value = 2.833
$ws.Cells.Item(13, 27).Value = 3.901
$ws.Cells.Item(13, 28).Value = 1.975
$ws.Cells.Item(13, 29).Value = 6.273
$ws.Cells.Item(13, 30).Value = 2.88
$ws.Cells.Item(13, 31).Value = 1.697
$ws.Cells.Item(13, 32).Value = 1.667
$ws.Cells.Item(13, 33).Value = 0.8
$ws.Cells.Item(13, 34).Value = 0.894
$ws.Cells.Item(13, 35).Value = 1.311
$ws.Cells.Item(13, 36).Value = 0.348
$ws.Cells.Item(13, 37).Value = 0.59
$ws.Cells.Item(13, 38).Value = 1.689
$ws.Cells.Item(13, 39).Value = 0.792
$ws.Cells.Item(13, 40).Value = 0.89
$ws.Cells.Item(13, 41).Value = 1.556
